# Auto-generated script to apply cryptos list price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "47.272.40"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -1.51%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.489.32"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.01%  "

$ws.Range("E4").Value = "  +0.05%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "321.66"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.88%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "108.46"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +2.10%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.522"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.81%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("E10").Value = "  +3.18%  "

$ws.Range("E11").Value = "  -1.02%  "

$ws.Range("E12").Value = "  +0.62%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "18.59"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.67%  "

$ws.Range("E14").Value = "  -0.27%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.877.81"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -1.09%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "2.495.34"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.846"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.31%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "47.222.79"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.31%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.47"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +5.18%  "

$ws.Range("E20").Value = "  +0.66%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0940"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.19%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "2.76"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +15.07%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "70.55"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.50%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "246.99"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.87%  "

$ws.Range("E25").Value = "  -0.93%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.09%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "25.73"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -2.54%  "

$ws.Range("E28").Value = "  +4.19%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "9.96"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.95%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.138"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +2.00%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "34.69"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.66%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "49.86"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.82%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "20.42"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +1.26%  "

$ws.Range("E34").Value = "  -1.13%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.0782"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.43%  "

$ws.Range("E36").Value = "  +0.15%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "4.80"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +2.64%  "

$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("E39").Value = "  -1.88%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "22.90"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +7.40%  "

$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("E42").Value = "  -2.00%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "119.76"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -1.31%  "

$ws.Range("E44").Value = "  -0.75%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.000.20"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +1.41%  "

$ws.Range("E46").Value = "  +0.58%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.00"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -4.25%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "9.14"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -1.07%  "

$ws.Range("E49").Value = "  -2.36%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "5.23"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -3.18%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "56.89"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +3.16%  "
